$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.932.52'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.035.78'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.13'
$ws.Range('D5').ClearFormats()
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.47'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +3.40%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.379'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0820'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '2.338.05'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.54'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.37'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.762'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '2.045.56'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '37.838.04'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.71'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('E20').Value = '  -6.81%  '
$ws.Range('D21').Value = '0.0₃0826'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.04'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.43'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.31'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.61'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.88'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('E30').Value = '  -3.16%  '
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.25'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +8.57%  '
$ws.Range('E33').Value = '  -2.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0606'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.51'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.34'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('E38').Value = '  +0.99%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.66'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.03%  '
$ws.Range('D41').Value = '1.533.86'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0218'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '96.42'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  -3.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0915'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.98'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.08'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = '2.226.96'
$ws.Range('E51').Value = '  -0.62%  '
